$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price (D) column cells being updated so Excel
# does not reinterpret numeric-looking strings as floating point numbers.
$priceCells = @("D2","D3","D5","D6","D7","D9","D10","D11","D13","D14","D15","D16","D17","D18","D20","D21","D22","D23","D24","D25","D27","D29","D30","D31","D32","D33","D34","D35","D39","D43","D47","D48","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '51.888.64'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '2.944.25'
$ws.Range('E3').Value = '  +4.14%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '353.17'
$ws.Range('E5').Value = '  +0.93%  '
$ws.Range('D6').Value = '112.65'
$ws.Range('E6').Value = '  -0.23%  '
$ws.Range('D7').Value = '0.562'
$ws.Range('E7').Value = '  +0.86%  '
$ws.Range('D9').Value = '0.629'
$ws.Range('E9').Value = '  +1.64%  '
$ws.Range('D10').Value = '39.46'
$ws.Range('E10').Value = '  -1.76%  '
$ws.Range('D11').Value = '0.0893'
$ws.Range('E11').Value = '  +5.49%  '
$ws.Range('E12').Value = '  +1.31%  '
$ws.Range('D13').Value = '19.98'
$ws.Range('E13').Value = '  -0.62%  '
$ws.Range('D14').Value = '7.88'
$ws.Range('E14').Value = '  +1.31%  '
$ws.Range('D15').Value = '3.405.41'
$ws.Range('E15').Value = '  +4.04%  '
$ws.Range('D16').Value = '2.938.09'
$ws.Range('E16').Value = '  +3.12%  '
$ws.Range('D17').Value = '0.992'
$ws.Range('E17').Value = '  +1.21%  '
$ws.Range('D18').Value = '51.962.15'
$ws.Range('E18').Value = '  +0.30%  '
$ws.Range('E19').Value = '  +0.67%  '
$ws.Range('D20').Value = '3.33'
$ws.Range('E20').Value = '  -2.93%  '
$ws.Range('D21').Value = '14.47'
$ws.Range('E21').Value = '  +7.71%  '
$ws.Range('D22').Value = '0.0₃0991'
$ws.Range('E22').Value = '  +1.97%  '
$ws.Range('D23').Value = '71.31'
$ws.Range('E23').Value = '  +1.28%  '
$ws.Range('D24').Value = '270.20'
$ws.Range('E24').Value = '  +0.52%  '
$ws.Range('D25').Value = '2.80'
$ws.Range('E25').Value = '  +1.99%  '
$ws.Range('E26').Value = '  +8.96%  '
$ws.Range('D27').Value = '26.98'
$ws.Range('E27').Value = '  +2.84%  '
$ws.Range('E28').Value = '  +0.13%  '
$ws.Range('D29').Value = '7.37'
$ws.Range('E29').Value = '  +16.54%  '
$ws.Range('D30').Value = '0.108'
$ws.Range('E30').Value = '  +20.39%  '
$ws.Range('D31').Value = '10.65'
$ws.Range('E31').Value = '  +1.21%  '
$ws.Range('D32').Value = '37.55'
$ws.Range('E32').Value = '  -3.44%  '
$ws.Range('D33').Value = '2.28'
$ws.Range('E33').Value = '  +0.69%  '
$ws.Range('D34').Value = '6.21'
$ws.Range('E34').Value = '  +10.18%  '
$ws.Range('D35').Value = '52.89'
$ws.Range('E35').Value = '  +0.12%  '
$ws.Range('E36').Value = '  +0.47%  '
$ws.Range('E37').Value = '  -0.16%  '
$ws.Range('E38').Value = '  +3.32%  '
$ws.Range('D39').Value = '18.85'
$ws.Range('E39').Value = '  -0.12%  '
$ws.Range('E40').Value = '  +1.83%  '
$ws.Range('E41').Value = '  +6.69%  '
$ws.Range('E42').Value = '  +1.86%  '
$ws.Range('D43').Value = '23.47'
$ws.Range('E43').Value = '  +5.85%  '
$ws.Range('E44').Value = '  -0.91%  '
$ws.Range('E45').Value = '  +0.11%  '
$ws.Range('E46').Value = '  +1.26%  '
$ws.Range('D47').Value = '2.170.11'
$ws.Range('D48').Value = '112.19'
$ws.Range('E48').Value = '  -8.67%  '
$ws.Range('E49').Value = '  -0.56%  '
$ws.Range('D50').Value = '0.0343'
$ws.Range('E50').Value = '  +10.82%  '
$ws.Range('D51').Value = '0.938'
$ws.Range('E51').Value = '  -0.51%  '
